# Insert a new data row at row 141, shifting existing rows 141:203 down to 142:204.
# The new row 141 gets fresh data; all rows that were 141..203 become 142..204 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 141 (pushes old row 141 and below down by one row).
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with its data. Columns A,B,C,E,F,G,H,N,O,Q,R
# are identical to the surrounding rows (constant for this market/category), and
# Excel's Insert() already copies formatting (incl. the date style) from the row above.
$ws.Cells.Item(141, 1).Value = 3
$ws.Cells.Item(141, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(141, 3).Value = "Coquimbo"
$ws.Cells.Item(141, 4).Value = 44489
$ws.Cells.Item(141, 5).Value = 5
$ws.Cells.Item(141, 6).Value = 100112012
$ws.Cells.Item(141, 7).Value = "Espinaca"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 120
$ws.Cells.Item(141, 11).Value = 2500
$ws.Cells.Item(141, 12).Value = 2500
$ws.Cells.Item(141, 13).Value = 2500
$ws.Cells.Item(141, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(141, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(141, 16).Value = 833
$ws.Cells.Item(141, 17).Value = 3
$ws.Cells.Item(141, 18).Value = "Hortaliza"
